$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Text = "4.53% and `$20.76"
$find.Replacement.ClearFormatting()
$find.Replacement.Text = "4.53%, and `$20.76"
$find.Execute($find.Text, $true, $false, $false, $false, $false, $true, 1, $false, $find.Replacement.Text, 2)
